$wb = $excel.ActiveWorkbook

# --- workbook-level window geometry ---
$wb.Windows.Item(1).Left = 4290
$wb.Windows.Item(1).Top = 4290
$wb.Windows.Item(1).Width = 28800

# --- notifications sheet: drop the duplicate "COURSE ENROLLMENT" rows ---
$wsNotif = $wb.Worksheets.Item("notifications")
$wsNotif.Rows.Item(2).Resize(2).Delete() | Out-Null
$wsNotif.Range("H1").Value = 0

# --- student_courses sheet: drop the duplicate enrollment row, update count ---
$wsCourses = $wb.Worksheets.Item("student_courses")
$wsCourses.Range("A3:E3").ClearContents() | Out-Null
$wsCourses.Range("L4").Value = 1
$wsCourses.Range("F17").Select() | Out-Null
